$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-based, including header row 1) that correspond to movies being
# removed from the list: Academy Award Review of Walt Disney Cartoons,
# Fantasia, Saludos Amigos, Victory Through Air Power, The Three Caballeros,
# Make Mine Music, Song of the South, Fun and Fancy Free, Melody Time,
# The Adventures of Ichabod and Mr. Toad, Fantasia 2000.
$rowsToDelete = @(2, 5, 8, 9, 10, 11, 12, 13, 14, 15, 43)

# Delete from the bottom up so earlier (lower-numbered) row indices stay valid.
$sorted = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sorted) {
    $ws.Rows($r).Delete()
}

# Re-apply the AutoFilter over the shrunk data range (header + 59 data rows).
$ws.AutoFilterMode = $false
$ws.Range("A1:E60").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase defined name in sync with
# the new AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet3!_FilterDatabase") {
        $n.RefersTo = "=Sheet3!`$A`$1:`$E`$60"
    }
}

# Move the active selection to A2, matching the author's final cursor spot.
$ws.Range("A2").Select()
